# ftest.xlsx - "added tests for calcrule 7 and 8"
#
# 1. Row 24 (fm19) description text is revised.
# 2. Row 33 (fm28) Windows/Linux status flips from "in progress" to "complete".
# 3. Two brand-new rows (34 = fm29, 35 = fm30) are appended to the ftests table,
#    formatted the same as the existing data rows.
# 4. The active selection moves down to the newly-added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- 1. fm19 description edit (row 24) ---------------------------------
$ws.Range("C24").Value = "Residential policy with blanket policy terms. Prior level loss back-allocation"

# --- 2. fm28 status edit (row 33) ---------------------------------------
$ws.Range("H33").Value = "complete"
$ws.Range("I33").Value = "complete"

# --- 3. New rows 34 & 35 -------------------------------------------------
# Clone the formatting of the last existing data row (32) onto the two new
# rows, then fill in the new test case data.
$ws.Range("B32:I32").Copy()
$ws.Range("B34:I35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B34").Value = "fm29"
$ws.Range("C34").Value = "Location deductibles with overall maximum policy deductible, and policy limit using calcrules 6 & 7"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = "6,7"
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = "in progress"
$ws.Range("I34").Value = "in progress"

$ws.Range("B35").Value = "fm30"
$ws.Range("C35").Value = "Location deductibles with overall minimum policy deductible, and policy limit using calcrule 6 & 8"
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = "6,8"
$ws.Range("F35").Value = 2
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = "in progress"
$ws.Range("I35").Value = "in progress"

# --- 4. Move the selection to the new last row ---------------------------
$ws.Activate()
$ws.Range("B35").Select()
